# Update countries & provincias Spain
# Re-sorted rankings (by total cases) shifted several countries' rows,
# and refreshed the day's case/death figures for the affected rows.
# Also bumps the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 14:34"

# --- Helper: write a full data row (country name + 7 numeric columns) ---
function Set-CountryRow($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Row 19 - Paises Bajos (figures refreshed)
Set-CountryRow 19 "Paises Bajos" 41774 455 0 36236 628 84 5288

# Row 24 - Portugal (figures refreshed)
Set-CountryRow 24 "Portugal" 26715 533 2258 23352 135 16 1105

# Rows 25-26 - Suecia now ranks above Pakistan
Set-CountryRow 25 "Suecia" 24623 705 4074 17509 425 99 3040
Set-CountryRow 26 "Pakistan" 24073 859 6464 17045 111 20 564

# Row 30 - Bielorrusia (figures refreshed)
Set-CountryRow 30 "Bielorrusia" 20168 913 5067 14985 92 4 116

# Rows 53-54 - Kuwait now ranks above Malasia
Set-CountryRow 53 "Kuwait" 6567 278 2381 4142 91 2 44
Set-CountryRow 54 "Malasia" 6467 39 4776 1584 19 0 107

# Row 75 - Croacia (figures refreshed)
Set-CountryRow 75 "Croacia" 2125 6 1641 398 14 1 86

# Rows 85-89 - Senegal now ranks above Nueva Zelanda, Honduras, Eslovenia, Eslovaquia
Set-CountryRow 85 "Senegal" 1492 59 562 917 6 1 13
Set-CountryRow 86 "Nueva Zelanda" 1489 1 1332 136 2 0 21
Set-CountryRow 87 "Honduras" 1461 191 132 1230 10 6 99
Set-CountryRow 88 "Eslovenia" 1449 1 247 1103 13 0 99
Set-CountryRow 89 "Eslovaquia" 1445 16 806 613 4 1 26

# Rows 101-102 - Sri Lanka now ranks above Guatemala
Set-CountryRow 101 "Sri Lanka" 804 7 232 563 1 0 9
Set-CountryRow 102 "Guatemala" 798 35 86 691 5 2 21

# Row 136 - Sierra Leona (figures refreshed)
Set-CountryRow 136 "Sierra Leona" 231 6 54 161 0 2 16
